# Scheduled-runner price refresh for the Leve profit sheets.
# For each affected Leve row, currentAveragePrice* (H/I/J) and the derived
# LevePrice*/LeveProfit* columns (K/L/M/N) are refreshed with newly fetched
# market-board figures. Column layout (all 8 sheets share the same header):
#   H = currentAveragePrice     K = LevePriceNQ
#   I = currentAveragePriceNQ   L = LevePriceHQ
#   J = currentAveragePriceHQ   M = LeveProfitNQ
#                                N = LeveProfitHQ

$wb = $excel.ActiveWorkbook

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 12: Don't Be So Tallow | Beeswax
$ws.Cells.Item(12, 8).Value = 1087.1428
$ws.Cells.Item(12, 10).Value = 208.33333
$ws.Cells.Item(12, 12).Value = 208.33333
$ws.Cells.Item(12, 14).Value = -548.3333299999999

# Row 33: Glazed and Confused | Clear Glass Lens
$ws.Cells.Item(33, 8).Value = 173.44444
$ws.Cells.Item(33, 9).Value = 176.92
$ws.Cells.Item(33, 11).Value = 176.92
$ws.Cells.Item(33, 13).Value = 52.08000000000001

# Row 43: Growing Is Knowing | Growth Formula Gamma
$ws.Cells.Item(43, 8).Value = 2794.5
$ws.Cells.Item(43, 9).Value = 559.375
$ws.Cells.Item(43, 10).Value = 3912.0625
$ws.Cells.Item(43, 11).Value = 559.375
$ws.Cells.Item(43, 12).Value = 3912.0625
$ws.Cells.Item(43, 13).Value = -490.375
$ws.Cells.Item(43, 14).Value = -4050.0625

# Row 70: Consecrating Congregation | Holy Water
$ws.Cells.Item(70, 8).Value = 1393.3334
$ws.Cells.Item(70, 10).Value = 1090
$ws.Cells.Item(70, 12).Value = 3270
$ws.Cells.Item(70, 14).Value = -3810

# Row 73: Curbing the Contagion (L) | Holy Water
$ws.Cells.Item(73, 8).Value = 1393.3334
$ws.Cells.Item(73, 10).Value = 1090
$ws.Cells.Item(73, 12).Value = 3270
$ws.Cells.Item(73, 14).Value = -5142

# Row 87: There Was a Late Fee | Noble Gold
$ws.Cells.Item(87, 8).Value = 39543.2
$ws.Cells.Item(87, 10).Value = 39543.2
$ws.Cells.Item(87, 12).Value = 39543.2
$ws.Cells.Item(87, 14).Value = -42039.2

# Row 90: A Gate Arcane Is Dragon's Bane (L) | Noble Gold
$ws.Cells.Item(90, 8).Value = 39543.2
$ws.Cells.Item(90, 10).Value = 39543.2
$ws.Cells.Item(90, 12).Value = 118629.6
$ws.Cells.Item(90, 14).Value = -131109.6

# Row 100: Asking for a Friend | Beetle Glue
$ws.Cells.Item(100, 8).Value = 2233.111
$ws.Cells.Item(100, 9).Value = 1624.5
$ws.Cells.Item(100, 10).Value = 2720
$ws.Cells.Item(100, 11).Value = 1624.5
$ws.Cells.Item(100, 12).Value = 2720
$ws.Cells.Item(100, 13).Value = -1083.5
$ws.Cells.Item(100, 14).Value = -3802

# Row 116: Growing Up | Growth Formula Kappa
$ws.Cells.Item(116, 8).Value = 27783612
$ws.Cells.Item(116, 9).Value = 125001750
$ws.Cells.Item(116, 10).Value = 7000.857
$ws.Cells.Item(116, 11).Value = 125001750
$ws.Cells.Item(116, 12).Value = 7000.857
$ws.Cells.Item(116, 13).Value = -124998308
$ws.Cells.Item(116, 14).Value = -13884.857

# Row 129: Practical Command | Commanding Craftsman's Draught
$ws.Cells.Item(129, 8).Value = 1254.9524
$ws.Cells.Item(129, 10).Value = 1275.8049
$ws.Cells.Item(129, 12).Value = 3827.4147
$ws.Cells.Item(129, 14).Value = -13827.4147

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 143313.61
$ws.Cells.Item(138, 10).Value = 160841.1
$ws.Cells.Item(138, 12).Value = 482523.3
$ws.Cells.Item(138, 14).Value = -492803.3

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 2: Ain't Got No Ingots | Bronze Ingot
$ws.Cells.Item(2, 8).Value = 1349.2307
$ws.Cells.Item(2, 9).Value = 1437.7778
$ws.Cells.Item(2, 11).Value = 1437.7778
$ws.Cells.Item(2, 13).Value = -1324.7778

# Row 32: Ingot We Trust | Steel Ingot
$ws.Cells.Item(32, 8).Value = 9778.918
$ws.Cells.Item(32, 9).Value = 6872.1206
$ws.Cells.Item(32, 10).Value = 21018.533
$ws.Cells.Item(32, 11).Value = 6872.1206
$ws.Cells.Item(32, 12).Value = 21018.533
$ws.Cells.Item(32, 13).Value = -6585.1206
$ws.Cells.Item(32, 14).Value = -21592.533

# Row 116: No Scope | Titanbronze Ingot
$ws.Cells.Item(116, 8).Value = 1349.2307
$ws.Cells.Item(116, 9).Value = 1437.7778
$ws.Cells.Item(116, 11).Value = 1437.7778
$ws.Cells.Item(116, 13).Value = 856.2221999999999

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 16442.611
$ws.Cells.Item(132, 9).Value = 2196.2144
$ws.Cells.Item(132, 11).Value = 6588.6432
$ws.Cells.Item(132, 13).Value = -4058.6432

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 3: Hells Bells | Bronze Ingot
$ws.Cells.Item(3, 8).Value = 1349.2307
$ws.Cells.Item(3, 9).Value = 1437.7778
$ws.Cells.Item(3, 11).Value = 1437.7778
$ws.Cells.Item(3, 13).Value = -1323.7778

# Row 7: Thank You for Your Business | Bronze Bastard Sword
$ws.Cells.Item(7, 8).Value = 3333566.2
$ws.Cells.Item(7, 9).Value = 3333566.2
$ws.Cells.Item(7, 11).Value = 3333566.2
$ws.Cells.Item(7, 13).Value = -3333453.2

# Row 43: Don't Fear the Reaper | Steel Scythe
$ws.Cells.Item(43, 8).Value = 119995
$ws.Cells.Item(43, 10).Value = 119995
$ws.Cells.Item(43, 12).Value = 119995
$ws.Cells.Item(43, 14).Value = -120357

# Row 99: Meddle in Metal | Oroshigane Ingot
$ws.Cells.Item(99, 8).Value = 1445.9166
$ws.Cells.Item(99, 9).Value = 1307
$ws.Cells.Item(99, 11).Value = 1307
$ws.Cells.Item(99, 13).Value = 191

# Row 105: Ingot to Wing It | Molybdenum Ingot
$ws.Cells.Item(105, 8).Value = 1043654.2
$ws.Cells.Item(105, 9).Value = 1616.3158
$ws.Cells.Item(105, 11).Value = 1616.3158
$ws.Cells.Item(105, 13).Value = 130.6841999999999

# Row 138: Bladewinner | Titanium Gold Greatsword
$ws.Cells.Item(138, 8).Value = 50673.6
$ws.Cells.Item(138, 10).Value = 50673.6
$ws.Cells.Item(138, 12).Value = 50673.6
$ws.Cells.Item(138, 14).Value = -60953.6

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 58: You Do the Heavy Lifting | Mahogany Lumber
$ws.Cells.Item(58, 8).Value = 23203.174
$ws.Cells.Item(58, 10).Value = 73003.86
$ws.Cells.Item(58, 12).Value = 73003.86
$ws.Cells.Item(58, 14).Value = -73409.86

# Row 136: Turali Quality | Dark Mahogany Lumber
$ws.Cells.Item(136, 8).Value = 23203.174
$ws.Cells.Item(136, 10).Value = 73003.86
$ws.Cells.Item(136, 12).Value = 219011.58
$ws.Cells.Item(136, 14).Value = -224111.58

# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Cells.Item(141, 8).Value = 28097.38
$ws.Cells.Item(141, 10).Value = 28097.38
$ws.Cells.Item(141, 12).Value = 28097.38
$ws.Cells.Item(141, 14).Value = -38457.38

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 619.1900000000001
$ws.Cells.Item(131, 10).Value = 748.9589
$ws.Cells.Item(131, 12).Value = 2246.8767
$ws.Cells.Item(131, 14).Value = -12326.8767

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 126: Gold Rush Order | Phrygian Gold Ingot
$ws.Cells.Item(126, 8).Value = 5385.5713
$ws.Cells.Item(126, 9).Value = 5006.0625
$ws.Cells.Item(126, 10).Value = 6600
$ws.Cells.Item(126, 11).Value = 15018.1875
$ws.Cells.Item(126, 12).Value = 19800
$ws.Cells.Item(126, 13).Value = -12548.1875
$ws.Cells.Item(126, 14).Value = -24740

# Row 132: On Board for Lar | Lar Ingot
$ws.Cells.Item(132, 8).Value = 19763.188
$ws.Cells.Item(132, 9).Value = 5730.6665
$ws.Cells.Item(132, 10).Value = 32144.824
$ws.Cells.Item(132, 11).Value = 17191.9995
$ws.Cells.Item(132, 12).Value = 96434.47200000001
$ws.Cells.Item(132, 13).Value = -14661.9995
$ws.Cells.Item(132, 14).Value = -101494.472

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 93: Hide to Go Seek | Gagana Leather
$ws.Cells.Item(93, 8).Value = 2122.4211
$ws.Cells.Item(93, 9).Value = 2029.2
$ws.Cells.Item(93, 10).Value = 2472
$ws.Cells.Item(93, 11).Value = 2029.2
$ws.Cells.Item(93, 12).Value = 2472
$ws.Cells.Item(93, 13).Value = -781.2
$ws.Cells.Item(93, 14).Value = -4968

# Row 100: Tiger in the Sack | Tiger Leather
$ws.Cells.Item(100, 8).Value = 2093.4666
$ws.Cells.Item(100, 9).Value = 1500.8334
$ws.Cells.Item(100, 10).Value = 2488.5557
$ws.Cells.Item(100, 11).Value = 1500.8334
$ws.Cells.Item(100, 12).Value = 2488.5557
$ws.Cells.Item(100, 13).Value = -959.8334
$ws.Cells.Item(100, 14).Value = -3570.5557

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 51: After the Smock-down | Linen Smock
# LeveProfitHQ (N) no longer applies -- cleared, LeveProfitNQ (M) added instead.
$ws.Cells.Item(51, 14).ClearContents() | Out-Null
$ws.Cells.Item(51, 8).Value = 8000
$ws.Cells.Item(51, 9).Value = 8000
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 8000
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -7490

# Row 81: Where the Dragonflies, the Net Catches | Crawler Silk
# LeveProfitHQ (N) cleared; LeveProfitNQ (M) keeps the refreshed figure.
$ws.Cells.Item(81, 14).ClearContents() | Out-Null
$ws.Cells.Item(81, 8).Value = 267.625
$ws.Cells.Item(81, 9).Value = 267.625
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 535.25
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = 525.75

# Row 84: To Kill a Dragon on Nameday (L) | Crawler Silk
# LeveProfitHQ (N) cleared; LeveProfitNQ (M) keeps the refreshed figure.
$ws.Cells.Item(84, 14).ClearContents() | Out-Null
$ws.Cells.Item(84, 8).Value = 267.625
$ws.Cells.Item(84, 9).Value = 267.625
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 2676.25
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 2627.75
